$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19 (ALC)
$ws.Range("H19").Value = 576.9091
$ws.Range("I19").Value = 616.3333
$ws.Range("J19").Value = 529.6
$ws.Range("K19").Value = 616.3333
$ws.Range("L19").Value = 529.6
$ws.Range("M19").Value = -441.3333
$ws.Range("N19").Value = -879.6

# Row 33 (ALC)
$ws.Range("H33").Value = 127.333336
$ws.Range("I33").Value = 112.8
$ws.Range("K33").Value = 112.8
$ws.Range("M33").Value = 116.2

# Row 106 (ALC)
$ws.Range("H106").Value = 17748.572
$ws.Range("I106").Value = 3936.25
$ws.Range("K106").Value = 3936.25
$ws.Range("M106").Value = -3305.25

# Row 113 (ALC)
$ws.Range("H113").Value = 66937
$ws.Range("I113").Value = 128249.125
$ws.Range("J113").Value = 5624.875
$ws.Range("K113").Value = 128249.125
$ws.Range("L113").Value = 5624.875
$ws.Range("M113").Value = -124995.125
$ws.Range("N113").Value = -12132.875

# Row 132 (ALC)
$ws.Range("H132").Value = 2389.2334
$ws.Range("I132").Value = 2389.2334
$ws.Range("K132").Value = 7167.7002
$ws.Range("M132").Value = -4637.7002

# Row 137 (ALC)
$ws.Range("H137").Value = 2260
$ws.Range("I137").Value = 2093.875
$ws.Range("J137").Value = 2555.3333
$ws.Range("K137").Value = 6281.625
$ws.Range("L137").Value = 7665.999899999999
$ws.Range("M137").Value = -3731.625
$ws.Range("N137").Value = -12765.9999

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (ARM)
$ws.Range("H2").Value = 3917.5
$ws.Range("I2").Value = 3502
$ws.Range("K2").Value = 3502
$ws.Range("M2").Value = -3389

# Row 32 (ARM)
$ws.Range("H32").Value = 20029.65
$ws.Range("I32").Value = 3471.3447
$ws.Range("J32").Value = 500220.5
$ws.Range("K32").Value = 3471.3447
$ws.Range("L32").Value = 500220.5
$ws.Range("M32").Value = -3184.3447
$ws.Range("N32").Value = -500794.5

# Row 45 (ARM)
$ws.Range("H45").Value = 6755.36
$ws.Range("I45").Value = 10430.833
$ws.Range("J45").Value = 3362.6155
$ws.Range("K45").Value = 10430.833
$ws.Range("L45").Value = 3362.6155
$ws.Range("M45").Value = -10053.833
$ws.Range("N45").Value = -4116.6155

# Row 61 (ARM)
$ws.Range("H61").Value = 2379.7368
$ws.Range("I61").Value = 2023.7407
$ws.Range("K61").Value = 2023.7407
$ws.Range("M61").Value = -1811.7407

# Row 63 (ARM)
$ws.Range("H63").Value = 3218
$ws.Range("I63").Value = 2979.1667
$ws.Range("J63").Value = 3576.25
$ws.Range("K63").Value = 2979.1667
$ws.Range("L63").Value = 3576.25
$ws.Range("M63").Value = -2293.1667
$ws.Range("N63").Value = -4948.25

# Row 66 (ARM)
$ws.Range("H66").Value = 3218
$ws.Range("I66").Value = 2979.1667
$ws.Range("J66").Value = 3576.25
$ws.Range("K66").Value = 14895.8335
$ws.Range("L66").Value = 17881.25
$ws.Range("M66").Value = -11463.8335
$ws.Range("N66").Value = -24745.25

# Row 116 (ARM)
$ws.Range("H116").Value = 3917.5
$ws.Range("I116").Value = 3502
$ws.Range("K116").Value = 3502
$ws.Range("M116").Value = -1208

# Row 124 (ARM)
$ws.Range("H124").Value = 79499.5
$ws.Range("J124").Value = 79499.5
$ws.Range("L124").Value = 79499.5
$ws.Range("N124").Value = -89319.5

# Row 132 (ARM)
$ws.Range("H132").Value = 2783.4194
$ws.Range("I132").Value = 1429
$ws.Range("K132").Value = 4287
$ws.Range("M132").Value = -1757

# Row 136 (ARM)
$ws.Range("H136").Value = 2379.7368
$ws.Range("I136").Value = 2023.7407
$ws.Range("K136").Value = 6071.2221
$ws.Range("M136").Value = -3521.2221

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (BSM)
$ws.Range("H3").Value = 3917.5
$ws.Range("I3").Value = 3502
$ws.Range("K3").Value = 3502
$ws.Range("M3").Value = -3388

# Row 9 (BSM)
$ws.Range("H9").Value = 29990
$ws.Range("J9").Value = 29990
$ws.Range("L9").Value = 29990
$ws.Range("N9").Value = -30326

# Row 24 (BSM)
$ws.Range("H24").Value = 10000
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()

# Row 34 (BSM)
$ws.Range("H34").Value = 10000
$ws.Range("J34").Value = 10000
$ws.Range("L34").Value = 10000
$ws.Range("N34").Value = -10228

# Row 86 (BSM)
$ws.Range("H86").Value = 1714.4348
$ws.Range("I86").Value = 1361.2
$ws.Range("K86").Value = 1361.2
$ws.Range("M86").Value = -238.2

# Row 89 (BSM)
$ws.Range("H89").Value = 1714.4348
$ws.Range("I89").Value = 1361.2
$ws.Range("K89").Value = 6806
$ws.Range("M89").Value = -1190

# Row 94 (BSM)
$ws.Range("H94").Value = 3967.9092
$ws.Range("I94").Value = 3726
$ws.Range("J94").Value = 5500
$ws.Range("K94").Value = 3726
$ws.Range("L94").Value = 5500
$ws.Range("M94").Value = -3275
$ws.Range("N94").Value = -6402

# Row 99 (BSM)
$ws.Range("H99").Value = 2089.0625
$ws.Range("I99").Value = 1952.4166
$ws.Range("J99").Value = 2499
$ws.Range("K99").Value = 1952.4166
$ws.Range("L99").Value = 2499
$ws.Range("M99").Value = -454.4166
$ws.Range("N99").Value = -5495

# Row 107 (BSM)
$ws.Range("H107").Value = 64508.938
$ws.Range("I107").Value = 92766.82000000001
$ws.Range("J107").Value = 2341.6
$ws.Range("K107").Value = 92766.82000000001
$ws.Range("L107").Value = 2341.6
$ws.Range("M107").Value = -90846.82000000001
$ws.Range("N107").Value = -6181.6

# Row 134 (BSM)
$ws.Range("H134").Value = 963.41174
$ws.Range("I134").Value = 825.2
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 2475.6
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = 59.39999999999964
$ws.Range("N134").Value = -11070

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 17349.695
$ws.Range("I31").Value = 3220
$ws.Range("J31").Value = 43842.875
$ws.Range("K31").Value = 3220
$ws.Range("L31").Value = 43842.875
$ws.Range("M31").Value = -2925
$ws.Range("N31").Value = -44432.875

# Row 34 (CRP)
$ws.Range("H34").Value = 17349.695
$ws.Range("I34").Value = 3220
$ws.Range("J34").Value = 43842.875
$ws.Range("K34").Value = 3220
$ws.Range("L34").Value = 43842.875
$ws.Range("M34").Value = -3018
$ws.Range("N34").Value = -44246.875

# Row 56 (CRP)
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()

# Row 58 (CRP)
$ws.Range("H58").Value = 1201.8
$ws.Range("I58").Value = 995
$ws.Range("K58").Value = 995
$ws.Range("M58").Value = -792

# Row 136 (CRP)
$ws.Range("H136").Value = 1201.8
$ws.Range("I136").Value = 995
$ws.Range("K136").Value = 2985
$ws.Range("M136").Value = -435

$ws = $wb.Worksheets.Item("CUL")
# Row 32 (CUL)
$ws.Range("H32").Value = 418.75
$ws.Range("I32").Value = 287.5
$ws.Range("J32").Value = 550
$ws.Range("K32").Value = 862.5
$ws.Range("L32").Value = 1650
$ws.Range("M32").Value = -579.5
$ws.Range("N32").Value = -2216

$ws = $wb.Worksheets.Item("GSM")
# Row 5 (GSM)
$ws.Range("H5").Value = 61.666668
$ws.Range("I5").Value = 61.666668
$ws.Range("K5").Value = 61.666668
$ws.Range("M5").Value = 50.333332

# Row 96 (GSM)
$ws.Range("H96").Value = 69927.664
$ws.Range("J96").Value = 69927.664
$ws.Range("L96").Value = 69927.664
$ws.Range("N96").Value = -75419.664

# Row 113 (GSM)
$ws.Range("H113").Value = 2312.5
$ws.Range("I113").Value = 2761.2
$ws.Range("J113").Value = 1992
$ws.Range("K113").Value = 2761.2
$ws.Range("L113").Value = 1992
$ws.Range("M113").Value = -591.1999999999998
$ws.Range("N113").Value = -6332

# Row 132 (GSM)
$ws.Range("H132").Value = 7464.6
$ws.Range("I132").Value = 6099.769
$ws.Range("K132").Value = 18299.307
$ws.Range("M132").Value = -15769.307

$ws = $wb.Worksheets.Item("LTW")
# Row 55 (LTW)
$ws.Range("H55").Value = 243
$ws.Range("J55").Value = 181.375
$ws.Range("L55").Value = 181.375
$ws.Range("N55").Value = -527.375

# Row 61 (LTW)
$ws.Range("H61").Value = 74793.66
$ws.Range("I61").Value = 74994.96000000001
$ws.Range("K61").Value = 74994.96000000001
$ws.Range("M61").Value = -74792.96000000001

# Row 100 (LTW)
$ws.Range("H100").Value = 15706.643
$ws.Range("I100").Value = 1582.5
$ws.Range("K100").Value = 1582.5
$ws.Range("M100").Value = -1041.5

# Row 113 (LTW)
$ws.Range("H113").Value = 74793.66
$ws.Range("I113").Value = 74994.96000000001
$ws.Range("K113").Value = 74994.96000000001
$ws.Range("M113").Value = -72824.96000000001

$ws = $wb.Worksheets.Item("WVR")
# Row 2 (WVR)
$ws.Range("H2").Value = 45750
$ws.Range("J2").Value = 50000
$ws.Range("L2").Value = 50000
$ws.Range("N2").Value = -50224

# Row 95 (WVR)
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

# Row 104 (WVR)
$ws.Range("H104").Value = 28185
$ws.Range("J104").Value = 28185
$ws.Range("L104").Value = 28185
$ws.Range("N104").Value = -35173

